# Update the "想去人数" (F column) figures across the workbook's sheets
# to reflect the refreshed data snapshot, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2980
$ws1.Range("F5").Value  = 77
$ws1.Range("F7").Value  = 2324
$ws1.Range("F8").Value  = 1684
$ws1.Range("F13").Value = 2664
$ws1.Range("F15").Value = 1529
$ws1.Range("F16").Value = 7082
$ws1.Range("F18").Value = 7228
$ws1.Range("F19").Value = 4
$ws1.Range("F20").Value = 11
$ws1.Range("F21").Value = 5482
$ws1.Range("F22").Value = 3113
$ws1.Range("F23").Value = 3486
$ws1.Range("F25").Value = 186
$ws1.Range("F26").Value = 1896
$ws1.Range("F28").Value = 301
$ws1.Range("F29").Value = 879
$ws1.Range("F31").Value = 284
$ws1.Range("F33").Value = 2425
$ws1.Range("F34").Value = 1201
$ws1.Range("F35").Value = 2732
$ws1.Range("F36").Value = 32
$ws1.Range("F39").Value = 392
$ws1.Range("F40").Value = 1085
$ws1.Range("F42").Value = 479
$ws1.Range("F43").Value = 526

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 216
$ws2.Range("F11").Value = 31
$ws2.Range("F12").Value = 146
$ws2.Range("F17").Value = 55

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 54

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 2980
$ws4.Range("F6").Value  = 77
$ws4.Range("F7").Value  = 2324
$ws4.Range("F8").Value  = 1684
$ws4.Range("F14").Value = 2664
$ws4.Range("F15").Value = 1529
$ws4.Range("F16").Value = 216
$ws4.Range("F18").Value = 31
$ws4.Range("F19").Value = 7082
$ws4.Range("F21").Value = 7228
$ws4.Range("F22").Value = 11
$ws4.Range("F23").Value = 5482
$ws4.Range("F24").Value = 3113
$ws4.Range("F25").Value = 3486
$ws4.Range("F29").Value = 1896
$ws4.Range("F31").Value = 55
$ws4.Range("F32").Value = 301
$ws4.Range("F33").Value = 879
$ws4.Range("F35").Value = 284
$ws4.Range("F37").Value = 2425
$ws4.Range("F38").Value = 1201
$ws4.Range("F40").Value = 2732
$ws4.Range("F41").Value = 32
$ws4.Range("F45").Value = 392
$ws4.Range("F46").Value = 1085
$ws4.Range("F48").Value = 479
$ws4.Range("F49").Value = 526
